# Auto-generated edit script: updates cryptos list cell values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.058.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.85%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.513.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.90%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''602.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.22%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''173.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +3.27%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = '''LidoStakedEther'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '''https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = '''3.515.12'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +1.18%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = '''XRP'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '''https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '''0.609'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.13%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.08%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.193'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -0.25%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''7.28'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +6.80%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.582'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +1.39%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''46.05'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.30%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''0.0000275'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -1.15%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''4.060.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +0.79%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''8.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.07%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''611.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -0.30%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''3.511.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.57%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''70.048.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +1.88%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  +0.87%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''17.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +0.55%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.873'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +0.28%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''8.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -19.06%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''15.58'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -1.35%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''96.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.29%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''3.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -2.12%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -0.01%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''2.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -1.71%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''34.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +4.05%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''8.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -1.59%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''8.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -3.34%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''2.99'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -2.89%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''651.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +14.05%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''6.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +1.34%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''1.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -3.28%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''3.61'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +3.34%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.0997'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -1.29%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''10.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +0.12%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.0478'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +9.46%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = '''OKB'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''56.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.84%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = '''FirstDigitalUSD'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''0.998'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.05%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  +3.82%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''3.326.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -1.85%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = '''PEPE'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = '''0.0₃0715'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +2.79%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = '''TheGraph'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = '''https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = '''0.311'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -3.56%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = '''ThetaToken'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''2.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +3.51%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = '''InjectiveProtocol'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''32.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.83%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''2.56'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +0.21%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.130'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +1.12%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''134.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +1.52%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +0.00%  '
$ws.Range("E51").Style = "Normal"
